$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H53").Value = 585.94446
$ws.Range("I53").Value = 391.18182
$ws.Range("J53").Value = 892
$ws.Range("K53").Value = 391.18182
$ws.Range("L53").Value = 892
$ws.Range("M53").Value = 245.81818
$ws.Range("N53").Value = -2166

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H100").Value = 2618.6155
$ws.Range("I100").Value = 2249.125
$ws.Range("J100").Value = 3209.8
$ws.Range("K100").Value = 2249.125
$ws.Range("L100").Value = 3209.8
$ws.Range("M100").Value = -1708.125
$ws.Range("N100").Value = -4291.8

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H101").Value = 777.0769
$ws.Range("I101").Value = 907.25
$ws.Range("J101").Value = 568.8
$ws.Range("K101").Value = 2721.75
$ws.Range("L101").Value = 1706.4
$ws.Range("M101").Value = -1099.75
$ws.Range("N101").Value = -4950.4

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H116").Value = 4367.684
$ws.Range("J116").Value = 4749.5
$ws.Range("L116").Value = 4749.5
$ws.Range("N116").Value = -11633.5

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H127").Value = 1499.9231
$ws.Range("I127").Value = 685
$ws.Range("J127").Value = 4216.3335
$ws.Range("K127").Value = 2055
$ws.Range("L127").Value = 12649.0005
$ws.Range("M127").Value = 2905
$ws.Range("N127").Value = -22569.0005

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H15").Value = 4637.6665
$ws.Range("J15").Value = 4637.6665
$ws.Range("L15").Value = 4637.6665
$ws.Range("N15").Value = -5337.6665

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 7564.973
$ws.Range("I32").Value = 8590.777
$ws.Range("K32").Value = 8590.777
$ws.Range("M32").Value = -8303.777

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H63").Value = 2341
$ws.Range("I63").Value = 2341
$ws.Range("K63").Value = 2341
$ws.Range("M63").Value = -1655

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H66").Value = 2341
$ws.Range("I66").Value = 2341
$ws.Range("K66").Value = 11705
$ws.Range("M66").Value = -8273

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H110").Value = 36326.1
$ws.Range("I110").Value = 44401.5
$ws.Range("K110").Value = 44401.5
$ws.Range("M110").Value = -42356.5

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H112").Value = 20225
$ws.Range("J112").Value = 20225
$ws.Range("L112").Value = 20225
$ws.Range("N112").Value = -23179

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H113").Value = 91193
$ws.Range("J113").Value = 91193
$ws.Range("L113").Value = 91193
$ws.Range("N113").Value = -99871

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H6").Value = 29994
$ws.Range("J6").Value = 29994
$ws.Range("L6").Value = 29994
$ws.Range("N6").Value = -30220

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 1784.9048
$ws.Range("I20").Value = 1655.0769
$ws.Range("J20").Value = 1995.875
$ws.Range("K20").Value = 1655.0769
$ws.Range("L20").Value = 1995.875
$ws.Range("M20").Value = -1408.0769
$ws.Range("N20").Value = -2489.875

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H22").Value = 525
$ws.Range("I22").Value = 525
$ws.Range("J22").Value = 0
$ws.Range("K22").Value = 525
$ws.Range("L22").Value = 0
$ws.Range("M22").Value = -352
$ws.Range("N22").ClearContents()

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H50").Value = 55418.5
$ws.Range("J50").Value = 55418.5
$ws.Range("L50").Value = 55418.5
$ws.Range("N50").Value = -56566.5

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H55").Value = 65286.668
$ws.Range("J55").Value = 65286.668
$ws.Range("L55").Value = 65286.668
$ws.Range("N55").Value = -65832.66800000001

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H99").Value = 3542.739
$ws.Range("I99").Value = 2169.1667
$ws.Range("J99").Value = 5041.1816
$ws.Range("K99").Value = 2169.1667
$ws.Range("L99").Value = 5041.1816
$ws.Range("M99").Value = -671.1667000000002
$ws.Range("N99").Value = -8037.1816

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H101").Value = 87332.336
$ws.Range("J101").Value = 87332.336
$ws.Range("L101").Value = 87332.336
$ws.Range("N101").Value = -93822.336

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H12").Value = 3996.3333
$ws.Range("J12").Value = 3994.5
$ws.Range("L12").Value = 3994.5
$ws.Range("N12").Value = -4334.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 1459.75
$ws.Range("I16").Value = 1473.6666
$ws.Range("J16").Value = 1445.8334
$ws.Range("K16").Value = 1473.6666
$ws.Range("L16").Value = 1445.8334
$ws.Range("M16").Value = -1186.6666
$ws.Range("N16").Value = -2019.8334

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 6170.029
$ws.Range("I31").Value = 3517.4707
$ws.Range("J31").Value = 8675.223
$ws.Range("K31").Value = 3517.4707
$ws.Range("L31").Value = 8675.223
$ws.Range("M31").Value = -3222.4707
$ws.Range("N31").Value = -9265.223

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H34").Value = 6170.029
$ws.Range("I34").Value = 3517.4707
$ws.Range("J34").Value = 8675.223
$ws.Range("K34").Value = 3517.4707
$ws.Range("L34").Value = 8675.223
$ws.Range("M34").Value = -3315.4707
$ws.Range("N34").Value = -9079.223

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H113").Value = 1459.75
$ws.Range("I113").Value = 1473.6666
$ws.Range("J113").Value = 1445.8334
$ws.Range("K113").Value = 1473.6666
$ws.Range("L113").Value = 1445.8334
$ws.Range("M113").Value = 696.3334
$ws.Range("N113").Value = -5785.8334

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H134").Value = 5225.5713
$ws.Range("I134").Value = 8000
$ws.Range("J134").Value = 4763.1665
$ws.Range("K134").Value = 24000
$ws.Range("L134").Value = 14289.4995
$ws.Range("M134").Value = -21465
$ws.Range("N134").Value = -19359.4995

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 2518.6875
$ws.Range("I5").Value = 1025.75
$ws.Range("K5").Value = 3077.25
$ws.Range("M5").Value = -2965.25

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H36").Value = 2044
$ws.Range("J36").Value = 3522
$ws.Range("L36").Value = 10566
$ws.Range("N36").Value = -10904

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H50").Value = 5017.5
$ws.Range("J50").Value = 8166.6665
$ws.Range("L50").Value = 24499.9995
$ws.Range("N50").Value = -25461.9995

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H53").Value = 5017.5
$ws.Range("J53").Value = 8166.6665
$ws.Range("L53").Value = 24499.9995
$ws.Range("N53").Value = -25461.9995

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H128").Value = 129999.336
$ws.Range("I128").Value = 129999.336
$ws.Range("K128").Value = 389998.008
$ws.Range("M128").Value = -385018.008

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H135").Value = 2518.6875
$ws.Range("I135").Value = 1025.75
$ws.Range("K135").Value = 9231.75
$ws.Range("M135").Value = -6696.75

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H45").Value = 40833.168
$ws.Range("J45").Value = 41250
$ws.Range("L45").Value = 41250
$ws.Range("N45").Value = -42368

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 5220.4
$ws.Range("I70").Value = 4844.75
$ws.Range("K70").Value = 4844.75
$ws.Range("M70").Value = -4574.75

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H73").Value = 5220.4
$ws.Range("I73").Value = 4844.75
$ws.Range("K73").Value = 4844.75
$ws.Range("M73").Value = -3908.75

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 1957.5
$ws.Range("I80").Value = 1959
$ws.Range("J80").Value = 1950
$ws.Range("K80").Value = 1959
$ws.Range("L80").Value = 1950
$ws.Range("M80").Value = -961
$ws.Range("N80").Value = -3946

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H83").Value = 1957.5
$ws.Range("I83").Value = 1959
$ws.Range("J83").Value = 1950
$ws.Range("K83").Value = 9795
$ws.Range("L83").Value = 1950
$ws.Range("M83").Value = -4803
$ws.Range("N83").Value = -19734

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 1304.3889
$ws.Range("I122").Value = 1305.6154
$ws.Range("J122").Value = 1301.2
$ws.Range("K122").Value = 3916.8462
$ws.Range("L122").Value = 3903.6
$ws.Range("M122").Value = -1466.8462
$ws.Range("N122").Value = -8803.6

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 1551.7142
$ws.Range("I132").Value = 1510.3334
$ws.Range("K132").Value = 4531.0002
$ws.Range("M132").Value = -2001.0002

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 3374.75
$ws.Range("I16").Value = 3599.75
$ws.Range("K16").Value = 3599.75
$ws.Range("M16").Value = -3429.75

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H61").Value = 3528.9546
$ws.Range("I61").Value = 3305.111
$ws.Range("J61").Value = 4536.25
$ws.Range("K61").Value = 3305.111
$ws.Range("L61").Value = 4536.25
$ws.Range("M61").Value = -3103.111
$ws.Range("N61").Value = -4940.25

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H82").Value = 2004.45
$ws.Range("I82").Value = 1092.6666
$ws.Range("K82").Value = 1092.6666
$ws.Range("M82").Value = -731.6666

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H85").Value = 2004.45
$ws.Range("I85").Value = 1092.6666
$ws.Range("K85").Value = 1092.6666
$ws.Range("M85").Value = 155.3334

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H98").Value = 102479.14
$ws.Range("J98").Value = 98670.8
$ws.Range("L98").Value = 98670.8
$ws.Range("N98").Value = -104660.8

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H110").Value = 20000
$ws.Range("J110").Value = 20000
$ws.Range("L110").Value = 20000
$ws.Range("N110").Value = -28180

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H113").Value = 3528.9546
$ws.Range("I113").Value = 3305.111
$ws.Range("J113").Value = 4536.25
$ws.Range("K113").Value = 3305.111
$ws.Range("L113").Value = 4536.25
$ws.Range("M113").Value = -1135.111
$ws.Range("N113").Value = -8876.25

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H98").Value = 0
$ws.Range("J98").Value = 0
$ws.Range("L98").Value = 0
$ws.Range("N98").ClearContents()
